$d = $word.ActiveDocument

# Locate the "Metalsmithing: Soldering Workshop, Quench Jewelry Arts,
# Minneapolis, MN." paragraph (end of the Residencies and Workshops
# section) and remove the first of the blank paragraphs that follow it,
# collapsing the extra spacing before "Professional Arts Experience".
$paras = $d.Paragraphs
$target = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "*Quench*Jewelry Arts*Minneapolis*MN.*") {
        $target = $i
        break
    }
}

if ($target -ne $null) {
    $blank = $paras.Item($target + 1)
    $blank.Range.Delete()
}
